# Apply the 17:13:39 scrape update (7 new rows on LP1912, 1 new row on LP1912-215,
# plus re-sort-driven row shuffles) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

# --- Sheet 1 data rows ---
$ws1.Cells.Item(2,1).Value2 = "Última actualización: 17:13:39"
$ws1.Cells.Item(3,1).Value2 = "Total filas: 279"
$ws1.Cells.Item(23,1).Value2 = "06:17:28"
$ws1.Cells.Item(23,2).Value2 = "07:21"
$ws1.Cells.Item(23,3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(23,4).Value2 = 64
$ws1.Cells.Item(23,5).Value2 = "LP1912"
$ws1.Cells.Item(24,1).Value2 = "05:57:13"
$ws1.Cells.Item(24,2).Value2 = "07:21"
$ws1.Cells.Item(24,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(24,4).Value2 = 84
$ws1.Cells.Item(24,5).Value2 = "LP1912"
$ws1.Cells.Item(40,1).Value2 = "06:35:22"
$ws1.Cells.Item(40,2).Value2 = "08:29"
$ws1.Cells.Item(40,3).Value2 = "15_ABASTO"
$ws1.Cells.Item(40,4).Value2 = 114
$ws1.Cells.Item(40,5).Value2 = "LP1912"
$ws1.Cells.Item(41,1).Value2 = "06:35:22"
$ws1.Cells.Item(41,2).Value2 = "08:29"
$ws1.Cells.Item(41,3).Value2 = "11_ETCHEVERRY"
$ws1.Cells.Item(41,4).Value2 = 114
$ws1.Cells.Item(41,5).Value2 = "LP1912"
$ws1.Cells.Item(58,1).Value2 = "08:10:18"
$ws1.Cells.Item(58,2).Value2 = "09:18"
$ws1.Cells.Item(58,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(58,4).Value2 = 68
$ws1.Cells.Item(58,5).Value2 = "LP1912"
$ws1.Cells.Item(59,1).Value2 = "07:38:39"
$ws1.Cells.Item(59,2).Value2 = "09:18"
$ws1.Cells.Item(59,3).Value2 = "15X38_ABASTO"
$ws1.Cells.Item(59,4).Value2 = 100
$ws1.Cells.Item(59,5).Value2 = "LP1912"
$ws1.Cells.Item(109,1).Value2 = "10:50:41"
$ws1.Cells.Item(109,2).Value2 = "11:53"
$ws1.Cells.Item(109,3).Value2 = "225_GOMEZ"
$ws1.Cells.Item(109,4).Value2 = 63
$ws1.Cells.Item(109,5).Value2 = "LP1912"
$ws1.Cells.Item(110,1).Value2 = "10:37:52"
$ws1.Cells.Item(110,2).Value2 = "11:53"
$ws1.Cells.Item(110,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(110,4).Value2 = 76
$ws1.Cells.Item(110,5).Value2 = "LP1912"
$ws1.Cells.Item(142,1).Value2 = "11:47:17"
$ws1.Cells.Item(142,2).Value2 = "12:48"
$ws1.Cells.Item(142,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(142,4).Value2 = 61
$ws1.Cells.Item(142,5).Value2 = "LP1912"
$ws1.Cells.Item(143,1).Value2 = "11:11:33"
$ws1.Cells.Item(143,2).Value2 = "12:48"
$ws1.Cells.Item(143,3).Value2 = "15X38_ABASTO"
$ws1.Cells.Item(143,4).Value2 = 97
$ws1.Cells.Item(143,5).Value2 = "LP1912"
$ws1.Cells.Item(241,1).Value2 = "15:17:33"
$ws1.Cells.Item(241,2).Value2 = "17:07"
$ws1.Cells.Item(241,3).Value2 = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(241,4).Value2 = 110
$ws1.Cells.Item(241,5).Value2 = "LP1912"
$ws1.Cells.Item(242,1).Value2 = "16:28:03"
$ws1.Cells.Item(242,2).Value2 = "17:07"
$ws1.Cells.Item(242,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(242,4).Value2 = 39
$ws1.Cells.Item(242,5).Value2 = "LP1912"
$ws1.Cells.Item(262,1).Value2 = "17:13:39"
$ws1.Cells.Item(262,2).Value2 = "17:51"
$ws1.Cells.Item(262,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(262,4).Value2 = 38
$ws1.Cells.Item(262,5).Value2 = "LP1912"
$ws1.Cells.Item(263,1).Value2 = "16:37:06"
$ws1.Cells.Item(263,2).Value2 = "17:52"
$ws1.Cells.Item(263,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(263,4).Value2 = 75
$ws1.Cells.Item(263,5).Value2 = "LP1912"
$ws1.Cells.Item(264,1).Value2 = "16:28:03"
$ws1.Cells.Item(264,2).Value2 = "17:53"
$ws1.Cells.Item(264,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(264,4).Value2 = 85
$ws1.Cells.Item(264,5).Value2 = "LP1912"
$ws1.Cells.Item(265,1).Value2 = "16:28:03"
$ws1.Cells.Item(265,2).Value2 = "17:58"
$ws1.Cells.Item(265,3).Value2 = "17_ROMERO"
$ws1.Cells.Item(265,4).Value2 = 90
$ws1.Cells.Item(265,5).Value2 = "LP1912"
$ws1.Cells.Item(266,1).Value2 = "16:14:44"
$ws1.Cells.Item(266,2).Value2 = "18:00"
$ws1.Cells.Item(266,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(266,4).Value2 = 106
$ws1.Cells.Item(266,5).Value2 = "LP1912"
$ws1.Cells.Item(267,1).Value2 = "16:14:44"
$ws1.Cells.Item(267,2).Value2 = "18:05"
$ws1.Cells.Item(267,3).Value2 = "11_ETCHEVERRY"
$ws1.Cells.Item(267,4).Value2 = 111
$ws1.Cells.Item(267,5).Value2 = "LP1912"
$ws1.Cells.Item(268,1).Value2 = "16:28:03"
$ws1.Cells.Item(268,2).Value2 = "18:06"
$ws1.Cells.Item(268,3).Value2 = "11_ETCHEVERRY"
$ws1.Cells.Item(268,4).Value2 = 98
$ws1.Cells.Item(268,5).Value2 = "LP1912"
$ws1.Cells.Item(270,1).Value2 = "16:14:44"
$ws1.Cells.Item(270,2).Value2 = "18:10"
$ws1.Cells.Item(270,3).Value2 = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(270,4).Value2 = 116
$ws1.Cells.Item(270,5).Value2 = "LP1912"
$ws1.Cells.Item(271,1).Value2 = "16:28:03"
$ws1.Cells.Item(271,2).Value2 = "18:17"
$ws1.Cells.Item(271,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(271,4).Value2 = 109
$ws1.Cells.Item(271,5).Value2 = "LP1912"
$ws1.Cells.Item(272,1).Value2 = "16:52:42"
$ws1.Cells.Item(272,2).Value2 = "18:21"
$ws1.Cells.Item(272,3).Value2 = "215C_EL PATO"
$ws1.Cells.Item(272,4).Value2 = 89
$ws1.Cells.Item(272,5).Value2 = "LP1912"
$ws1.Cells.Item(273,1).Value2 = "16:28:03"
$ws1.Cells.Item(273,2).Value2 = "18:22"
$ws1.Cells.Item(273,3).Value2 = "215C_EL PATO"
$ws1.Cells.Item(273,4).Value2 = 114
$ws1.Cells.Item(273,5).Value2 = "LP1912"
$ws1.Cells.Item(274,1).Value2 = "16:28:03"
$ws1.Cells.Item(274,2).Value2 = "18:25"
$ws1.Cells.Item(274,3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(274,4).Value2 = 117
$ws1.Cells.Item(274,5).Value2 = "LP1912"
$ws1.Cells.Item(275,1).Value2 = "16:52:42"
$ws1.Cells.Item(275,2).Value2 = "18:29"
$ws1.Cells.Item(275,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(275,4).Value2 = 97
$ws1.Cells.Item(275,5).Value2 = "LP1912"
$ws1.Cells.Item(276,1).Value2 = "16:37:06"
$ws1.Cells.Item(276,2).Value2 = "18:30"
$ws1.Cells.Item(276,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(276,4).Value2 = 113
$ws1.Cells.Item(276,5).Value2 = "LP1912"
$ws1.Cells.Item(277,1).Value2 = "16:37:06"
$ws1.Cells.Item(277,2).Value2 = "18:36"
$ws1.Cells.Item(277,3).Value2 = "15X38_ABASTO"
$ws1.Cells.Item(277,4).Value2 = 119
$ws1.Cells.Item(277,5).Value2 = "LP1912"
$ws1.Cells.Item(278,1).Value2 = "17:13:39"
$ws1.Cells.Item(278,2).Value2 = "18:36"
$ws1.Cells.Item(278,3).Value2 = "23_HERNANDEZ"
$ws1.Cells.Item(278,4).Value2 = 83
$ws1.Cells.Item(278,5).Value2 = "LP1912"
$ws1.Cells.Item(279,1).Value2 = "17:13:39"
$ws1.Cells.Item(279,2).Value2 = "18:41"
$ws1.Cells.Item(279,3).Value2 = "10_OLMOS"
$ws1.Cells.Item(279,4).Value2 = 88
$ws1.Cells.Item(279,5).Value2 = "LP1912"
$ws1.Cells.Item(280,1).Value2 = "16:52:42"
$ws1.Cells.Item(280,2).Value2 = "18:45"
$ws1.Cells.Item(280,3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(280,4).Value2 = 113
$ws1.Cells.Item(280,5).Value2 = "LP1912"
$ws1.Cells.Item(281,1).Value2 = "17:13:39"
$ws1.Cells.Item(281,2).Value2 = "18:52"
$ws1.Cells.Item(281,3).Value2 = "17_ROMERO"
$ws1.Cells.Item(281,4).Value2 = 99
$ws1.Cells.Item(281,5).Value2 = "LP1912"
$ws1.Cells.Item(282,1).Value2 = "17:13:39"
$ws1.Cells.Item(282,2).Value2 = "18:57"
$ws1.Cells.Item(282,3).Value2 = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(282,4).Value2 = 104
$ws1.Cells.Item(282,5).Value2 = "LP1912"
$ws1.Cells.Item(283,1).Value2 = "17:13:39"
$ws1.Cells.Item(283,2).Value2 = "18:59"
$ws1.Cells.Item(283,3).Value2 = "14_ABASTO"
$ws1.Cells.Item(283,4).Value2 = 106
$ws1.Cells.Item(283,5).Value2 = "LP1912"
$ws1.Cells.Item(284,1).Value2 = "17:13:39"
$ws1.Cells.Item(284,2).Value2 = "19:03"
$ws1.Cells.Item(284,3).Value2 = "215_EL PELIGRO"
$ws1.Cells.Item(284,4).Value2 = 110
$ws1.Cells.Item(284,5).Value2 = "LP1912"

# --- Sheet 2 data rows ---
$ws2.Cells.Item(2,1).Value2 = "Última actualización: 17:13:39"
$ws2.Cells.Item(3,1).Value2 = "Total filas: 44"
$ws2.Cells.Item(48,1).Value2 = "16:28:03"
$ws2.Cells.Item(48,2).Value2 = "18:22"
$ws2.Cells.Item(48,3).Value2 = "215C_EL PATO"
$ws2.Cells.Item(48,4).Value2 = 114
$ws2.Cells.Item(48,5).Value2 = "LP1912"
$ws2.Cells.Item(49,1).Value2 = "17:13:39"
$ws2.Cells.Item(49,2).Value2 = "19:03"
$ws2.Cells.Item(49,3).Value2 = "215_EL PELIGRO"
$ws2.Cells.Item(49,4).Value2 = 110
$ws2.Cells.Item(49,5).Value2 = "LP1912"

# --- Sheet 3 data rows ---
$ws3.Cells.Item(2,1).Value2 = "Última actualización: 17:13:39"

